$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.068.39"
Set-TextValue "E2" "  +2.15%  "
Set-TextValue "D3" "1.770.48"
Set-TextValue "E3" "  +0.83%  "
Set-TextValue "D4" "0.9970"
Set-TextValue "E4" "  -0.72%  "
Set-TextValue "D5" "322.39"
Set-TextValue "E5" "  -0.44%  "
Set-TextValue "D6" "0.9963"
Set-TextValue "E6" "  -0.52%  "
Set-TextValue "D7" "0.4254"
Set-TextValue "E7" "  -3.34%  "
Set-TextValue "D8" "0.3609"
Set-TextValue "E8" "  -2.61%  "
Set-TextValue "D9" "44.46"
Set-TextValue "E9" "  -0.93%  "
Set-TextValue "D10" "0.07484"
Set-TextValue "E10" "  -2.39%  "
Set-TextValue "D11" "1.112"
Set-TextValue "E11" "  +0.08%  "
Set-TextValue "D12" "0.9956"
Set-TextValue "E12" "  -0.65%  "
Set-TextValue "D13" "21.58"
Set-TextValue "E13" "  +0.08%  "
Set-TextValue "D14" "6.121"
Set-TextValue "E14" "  -0.47%  "
Set-TextValue "D15" "7.367"
Set-TextValue "E15" "  -0.56%  "
Set-TextValue "D16" "1.784.46"
Set-TextValue "E16" "  +1.21%  "
Set-TextValue "D17" "91.95"
Set-TextValue "E17" "  +1.93%  "
Set-TextValue "D18" "0.00001062"
Set-TextValue "E18" "  -0.88%  "
Set-TextValue "D19" "0.06389"
Set-TextValue "E19" "  +2.54%  "
Set-TextValue "D20" "0.9965"
Set-TextValue "E20" "  -0.47%  "
Set-TextValue "D21" "17.22"
Set-TextValue "E21" "  -0.93%  "
Set-TextValue "D22" "6.005"
Set-TextValue "E22" "  -2.68%  "
Set-TextValue "D23" "28.050.19"
Set-TextValue "E23" "  +1.93%  "
Set-TextValue "D24" "11.37"
Set-TextValue "E24" "  -1.12%  "
Set-TextValue "D25" "2.161"
Set-TextValue "E25" "  -6.06%  "
Set-TextValue "D26" "159.83"
Set-TextValue "E26" "  +4.67%  "
Set-TextValue "D27" "20.35"
Set-TextValue "E27" "  -0.55%  "
Set-TextValue "D28" "1.986.08"
Set-TextValue "E28" "  +1.34%  "
Set-TextValue "D29" "2.165"
Set-TextValue "E29" "  -5.33%  "
Set-TextValue "D30" "126.40"
Set-TextValue "E30" "  -0.65%  "
Set-TextValue "D31" "1.179"
Set-TextValue "E31" "  +0.64%  "
Set-TextValue "D32" "5.700"
Set-TextValue "E32" "  -0.01%  "
Set-TextValue "D33" "0.09030"
Set-TextValue "E33" "  -1.45%  "
Set-TextValue "D34" "3.502"
Set-TextValue "E34" "  -3.92%  "
Set-TextValue "D35" "12.69"
Set-TextValue "E35" "  +1.11%  "
Set-TextValue "D36" "0.02336"
Set-TextValue "E36" "  +1.20%  "
Set-TextValue "D37" "5.082"
Set-TextValue "E37" "  +0.66%  "
Set-TextValue "D38" "0.06094"
Set-TextValue "E38" "  -0.16%  "
Set-TextValue "D39" "0.2113"
Set-TextValue "E39" "  -1.61%  "
Set-TextValue "D40" "0.6428"
Set-TextValue "E40" "  +0.00%  "
Set-TextValue "D41" "1.189"
Set-TextValue "E41" "  +1.00%  "
Set-TextValue "D42" "0.9964"
Set-TextValue "E42" "  -0.44%  "
Set-TextValue "D45" "13.60"
Set-TextValue "E45" "  +0.24%  "
Set-TextValue "D46" "0.5989"
Set-TextValue "E46" "  +0.66%  "
Set-TextValue "D47" "3.703"
Set-TextValue "E47" "  -0.34%  "
Set-TextValue "D48" "2.005"
Set-TextValue "E48" "  +1.59%  "
Set-TextValue "D49" "123.86"
Set-TextValue "E49" "  -1.64%  "
Set-TextValue "D50" "1.184"
Set-TextValue "E50" "  +4.61%  "
Set-TextValue "D51" "0.06894"
Set-TextValue "E51" "  +0.34%  "

# Rows 43/44: coin name, link, price and volume all change (ranking swap)
Set-TextValue "B43" "FraxShare"
Set-TextValue "C43" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "7.886"
Set-TextValue "E43" "  -0.35%  "

Set-TextValue "B44" "WEMIXTOKEN"
Set-TextValue "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D44" "1.401"
Set-TextValue "E44" "  +1.21%  "
